$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 93, shifting rows 93:181 down to 94:182
$ws.Rows.Item(93).Insert()

# Populate the new row 93 with its data (same market/category as neighboring rows)
$ws.Range("A93").Value = 8
$ws.Range("B93").Value = "Terminal La Palmera de La Serena"
$ws.Range("C93").Value = "Coquimbo"
$ws.Range("D93").Value = 44673
$ws.Range("D93").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E93").Value = 4
$ws.Range("F93").Value = 100112037
$ws.Range("G93").Value = "Cebollín"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 1100
$ws.Range("L93").Value = 1200
$ws.Range("M93").Value = 1150
$ws.Range("N93").Value = "$/paquete 6 unidades"
$ws.Range("O93").Value = "Provincia del Elquí"
$ws.Range("P93").Value = 192
$ws.Range("Q93").Value = 6
$ws.Range("R93").Value = "Hortaliza"
